# Setup of color test cases:
# Insert a new Todo row into the "Active" sheet (row 10), pushing existing
# rows 10-41 down to 11-42, and bump the "Max Id" tracker on the Config sheet.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$config = $wb.Worksheets.Item("Config")

# Insert a brand new row at position 10 (shifts rows 10..41 down to 11..42)
$active.Rows.Item(10).Insert()

# Use a neighboring plain-text cell's style as a template so the new
# date-like text cell does not pick up an automatic date number format.
$plainStyle = $active.Range("C10").Style

$active.Range("A10").Value = 90
$active.Range("B10").Value = "name all Settings properties in OneImageForm with prefix ""Setting"" for consistency"
$active.Range("C10").Value = "Todo"
$active.Range("D10").Value = "Task"
$active.Range("E10").Value = "'8/25/2018"
$active.Range("E10").Style = $plainStyle

# Update the "Max Id" value tracked on the Config sheet.
$config.Range("F2").Value = 90

Write-Host "edit complete"
